$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Price" column: values are numeric-looking strings stored as text
#     (inlineStr) in the source workbook. Force Text format first so Excel
#     does not silently coerce these strings into real numbers. ---
$priceCells = @("D2","D3","D7","D8","D9","D10","D11","D12","D13","D14","D15","D16","D17","D19","D20","D21","D23","D26","D40","D41","D42","D43","D44","D45","D46","D47","D48","D49")
foreach ($cell in $priceCells) {
    $ws.Range($cell).NumberFormat = "@"
}

$ws.Range("D2").Value = "237.54"
$ws.Range("D3").Value = "21.54"
$ws.Range("D7").Value = "3.352"
$ws.Range("D8").Value = "1.073"
$ws.Range("D9").Value = "0.7926"
$ws.Range("D10").Value = "0.1399"
$ws.Range("D11").Value = "0.07349"
$ws.Range("D12").Value = "0.03191"
$ws.Range("D13").Value = "0.02965"
$ws.Range("D14").Value = "0.09255"
$ws.Range("D15").Value = "0.001659"
$ws.Range("D16").Value = "3.250"
$ws.Range("D17").Value = "0.04776"
$ws.Range("D19").Value = "0.006220"
$ws.Range("D20").Value = "0.005108"
$ws.Range("D21").Value = "0.001051"
$ws.Range("D23").Value = "3.885"
$ws.Range("D26").Value = "0.1056"
$ws.Range("D40").Value = "0.04116"
$ws.Range("D41").Value = "0.006936"
$ws.Range("D42").Value = "0.003502"
$ws.Range("D43").Value = "0.1040"
$ws.Range("D44").Value = "0.009804"
$ws.Range("D45").Value = "0.00005437"
$ws.Range("D46").Value = "0.00000000750"
$ws.Range("D47").Value = "0.6757"
$ws.Range("D48").Value = "0.03714"
$ws.Range("D49").Value = "0.00002101"

# --- Plain text columns (Coin name, Link, Volume(1h) label) ---
$ws.Range("E18").Value = "17OneONEWorstin24h"
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("E42").Value = "41CEJICEJI"
$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("E43").Value = "42BKEXTokenBKK"
$ws.Range("E48").Value = "47BOLOBOLO"
